$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEDDATA")
$r = $ws.Range("G14")
$r.Value = 1.4809504892378391
Write-Host "G14 after .Value= :" $r.Value2

$r2 = $ws.Range("G15")
$r2.Value2 = 1.4898155228409224
Write-Host "G15 after .Value2= :" $r2.Value2

$r3 = $ws.Range("K14")
Write-Host "K14 value:" $r3.Value2
Write-Host "K14 formula:" $r3.Formula
Write-Host "K14 HasError:" $r3.HasFormula
